$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($rowIndex, $text) {
    $cell = $tbl.Rows.Item($rowIndex).Cells.Item(1)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

# Simple single-value replacements (1-indexed rows)
Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "2101"
Set-CellText 5 "0.00001"
Set-CellText 6 "0.00055"
Set-CellText 7 "0.00017"
Set-CellText 9 "0.00033"
Set-CellText 10 "0.00034"
Set-CellText 11 "0.00036"
Set-CellText 12 "0.39611"

# Rows 44, 45, 46 (1-indexed) had multi-run tab-delimited content; replace with single value
Set-CellText 44 "99.88"
Set-CellText 45 "0.4"
Set-CellText 46 "335"
